# Apply the edit described by the diff:
#  - add three new columns (M: best, N: Worst, O: Avg) that summarize
#    the existing per-row metrics in columns B:K
#  - headers in row 1, values (max / min / rounded average) in rows 2-10
#  - formatting to match: bold header cells, centered/top aligned,
#    M column gets a left border, N/O columns share the bold/no-border look
#  - select M2:M10 as the active selection, matching the saved view state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataCol = 2   # column B
$lastDataCol  = 11  # column K
$firstRow     = 2
$lastRow      = 10

# ---- headers ---------------------------------------------------------
# Written in this order so new shared-strings entries land as
# 19=Avg, 20=Worst, 21=best (matching the target file's string table).
$ws.Cells.Item(1, 15).Value2 = "Avg"    # O1
$ws.Cells.Item(1, 14).Value2 = "Worst"  # N1
$ws.Cells.Item(1, 13).Value2 = "best"   # M1

# ---- per-row summary values -----------------------------------------
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $vals = @()
    for ($c = $firstDataCol; $c -le $lastDataCol; $c++) {
        $vals += $ws.Cells.Item($r, $c).Value2
    }

    $best  = ($vals | Measure-Object -Maximum).Maximum
    $worst = ($vals | Measure-Object -Minimum).Minimum
    $avg   = [Math]::Round((($vals | Measure-Object -Average).Average), 4)

    $ws.Cells.Item($r, 13).Value2 = $best   # M
    $ws.Cells.Item($r, 14).Value2 = $worst  # N
    $ws.Cells.Item($r, 15).Value2 = $avg    # O
}

# ---- formatting -------------------------------------------------------
# N1 ("Worst"): bold, centered horizontally, top-aligned, no border
$n1 = $ws.Range("N1")
$n1.Font.Bold = $true
$n1.HorizontalAlignment = -4108   # xlCenter
$n1.VerticalAlignment = -4160     # xlTop

# M1 ("best"): bold, centered horizontally, top-aligned, thin border on the left
$m1 = $ws.Range("M1")
$m1.Borders.Item(7).LineStyle = 1   # xlEdgeLeft, xlContinuous
$m1.Borders.Item(7).Weight = 2      # xlThin
$m1.Font.Bold = $true
$m1.HorizontalAlignment = -4108     # xlCenter
$m1.VerticalAlignment = -4160       # xlTop

# ---- selection, matching the workbook's saved view ---------------------
$ws.Range("M2:M10").Select()
